$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 797.3570999999999
$ws.Range("I33").Value = 312.57144
$ws.Range("K33").Value = 312.57144
$ws.Range("M33").Value = -83.57144

$ws.Range("H137").Value = 2443.4583
$ws.Range("I137").Value = 2551.2942
$ws.Range("J137").Value = 2181.5715
$ws.Range("K137").Value = 7653.882599999999
$ws.Range("L137").Value = 6544.7145
$ws.Range("M137").Value = -5103.882599999999
$ws.Range("N137").Value = -11644.7145

$ws.Range("H138").Value = 229029.6
$ws.Range("J138").Value = 419159.38
$ws.Range("L138").Value = 1257478.14
$ws.Range("N138").Value = -1267758.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3632.261
$ws.Range("I45").Value = 3030.1667
$ws.Range("K45").Value = 3030.1667
$ws.Range("M45").Value = -2653.1667

$ws.Range("H61").Value = 5673.1904
$ws.Range("I61").Value = 3695.5151
$ws.Range("J61").Value = 7848.6333
$ws.Range("K61").Value = 3695.5151
$ws.Range("L61").Value = 7848.6333
$ws.Range("M61").Value = -3483.5151
$ws.Range("N61").Value = -8272.633300000001

$ws.Range("H74").Value = 2224.9167
$ws.Range("I74").Value = 831.9474
$ws.Range("K74").Value = 831.9474
$ws.Range("M74").Value = 42.05259999999998

$ws.Range("H77").Value = 2224.9167
$ws.Range("I77").Value = 831.9474
$ws.Range("K77").Value = 4159.737
$ws.Range("M77").Value = 208.2629999999999

$ws.Range("H132").Value = 1564.3433
$ws.Range("J132").Value = 3914
$ws.Range("L132").Value = 11742
$ws.Range("N132").Value = -16802

$ws.Range("H136").Value = 5673.1904
$ws.Range("I136").Value = 3695.5151
$ws.Range("J136").Value = 7848.6333
$ws.Range("K136").Value = 11086.5453
$ws.Range("L136").Value = 23545.8999
$ws.Range("M136").Value = -8536.5453
$ws.Range("N136").Value = -28645.8999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2898.125
$ws.Range("I99").Value = 1970.8889
$ws.Range("J99").Value = 5679.8335
$ws.Range("K99").Value = 1970.8889
$ws.Range("L99").Value = 5679.8335
$ws.Range("M99").Value = -472.8888999999999
$ws.Range("N99").Value = -8675.833500000001

$ws.Range("H134").Value = 5662.7705
$ws.Range("I134").Value = 2601.175
$ws.Range("K134").Value = 7803.525000000001
$ws.Range("M134").Value = -5268.525000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3279.375
$ws.Range("I99").Value = 3036.8333
$ws.Range("K99").Value = 3036.8333
$ws.Range("M99").Value = -1538.8333

$ws.Range("H103").Value = 27724.2
$ws.Range("I103").Value = 27724.2
$ws.Range("K103").Value = 27724.2
$ws.Range("M103").Value = -26552.2

$ws.Range("H126").Value = 3279.375
$ws.Range("I126").Value = 3036.8333
$ws.Range("K126").Value = 9110.499899999999
$ws.Range("M126").Value = -6640.499899999999

$ws.Range("H132").Value = 626358.4
$ws.Range("I132").Value = 702977.8
$ws.Range("J132").Value = 2456.8572
$ws.Range("K132").Value = 2108933.4
$ws.Range("L132").Value = 7370.571599999999
$ws.Range("M132").Value = -2106403.4
$ws.Range("N132").Value = -12430.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 73229.92999999999
$ws.Range("I22").Value = 530.3333
$ws.Range("J22").Value = 204089.2
$ws.Range("K22").Value = 1590.9999
$ws.Range("L22").Value = 612267.6000000001
$ws.Range("M22").Value = -1421.9999
$ws.Range("N22").Value = -612605.6000000001

$ws.Range("H27").Value = 73229.92999999999
$ws.Range("I27").Value = 530.3333
$ws.Range("J27").Value = 204089.2
$ws.Range("K27").Value = 1590.9999
$ws.Range("L27").Value = 612267.6000000001
$ws.Range("M27").Value = -1488.9999
$ws.Range("N27").Value = -612471.6000000001

$ws.Range("H33").Value = 771.75
$ws.Range("I33").Value = 37
$ws.Range("J33").Value = 1016.6667
$ws.Range("K33").Value = 222
$ws.Range("L33").Value = 6100.0002
$ws.Range("M33").Value = 61
$ws.Range("N33").Value = -6666.0002

$ws.Range("H38").Value = 3529.0588
$ws.Range("J38").Value = 7117.375
$ws.Range("L38").Value = 21352.125
$ws.Range("N38").Value = -22046.125

$ws.Range("H58").Value = 12901
$ws.Range("J58").Value = 21166.666
$ws.Range("L58").Value = 63499.99800000001
$ws.Range("N58").Value = -63755.99800000001

$ws.Range("H92").Value = 840.3
$ws.Range("I92").Value = 404
$ws.Range("K92").Value = 1212
$ws.Range("M92").Value = 36

$ws.Range("H117").Value = 1773.75
$ws.Range("J117").Value = 3000
$ws.Range("L117").Value = 9000
$ws.Range("N117").Value = -15884

$ws.Range("H122").Value = 1193.0834
$ws.Range("J122").Value = 1165.2727
$ws.Range("L122").Value = 10487.4543
$ws.Range("N122").Value = -15387.4543

$ws.Range("H131").Value = 170956.36
$ws.Range("I131").Value = 667727.4399999999
$ws.Range("K131").Value = 2003182.32
$ws.Range("M131").Value = -1998142.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 76080
$ws.Range("I62").Value = 76075
$ws.Range("K62").Value = 76075
$ws.Range("M62").Value = -75389

$ws.Range("H65").Value = 76080
$ws.Range("I65").Value = 76075
$ws.Range("K65").Value = 228225
$ws.Range("M65").Value = -224793

$ws.Range("H113").Value = 2722.4
$ws.Range("I113").Value = 2149.75
$ws.Range("K113").Value = 2149.75
$ws.Range("M113").Value = 20.25

$ws.Range("H126").Value = 3108.1
$ws.Range("J126").Value = 2350
$ws.Range("L126").Value = 7050
$ws.Range("N126").Value = -11990

$ws.Range("H132").Value = 2432.373
$ws.Range("I132").Value = 2218.9795
$ws.Range("J132").Value = 3013.2778
$ws.Range("K132").Value = 6656.9385
$ws.Range("L132").Value = 9039.8334
$ws.Range("M132").Value = -4126.9385
$ws.Range("N132").Value = -14099.8334

$ws.Range("H140").Value = 97470
$ws.Range("I140").Value = 90000
$ws.Range("J140").Value = 99960
$ws.Range("K140").Value = 90000
$ws.Range("L140").Value = 99960
$ws.Range("N140").Value = -110320
$ws.Range("M140").Value = -84820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 724
$ws.Range("I16").Value = 471.6
$ws.Range("K16").Value = 471.6
$ws.Range("M16").Value = -301.6

$ws.Range("H22").Value = 2740.8845
$ws.Range("I22").Value = 644
$ws.Range("J22").Value = 3122.1365
$ws.Range("K22").Value = 644
$ws.Range("L22").Value = 3122.1365
$ws.Range("M22").Value = -349
$ws.Range("N22").Value = -3712.1365

$ws.Range("H27").Value = 2740.8845
$ws.Range("I27").Value = 644
$ws.Range("J27").Value = 3122.1365
$ws.Range("K27").Value = 644
$ws.Range("L27").Value = 3122.1365
$ws.Range("M27").Value = -537
$ws.Range("N27").Value = -3336.1365

$ws.Range("H40").Value = 3917.2563
$ws.Range("I40").Value = 3963.4644
$ws.Range("J40").Value = 3799.6365
$ws.Range("K40").Value = 3963.4644
$ws.Range("L40").Value = 3799.6365
$ws.Range("M40").Value = -3827.4644
$ws.Range("N40").Value = -4071.6365

$ws.Range("H122").Value = 2893.2974
$ws.Range("I122").Value = 2239.7917
$ws.Range("J122").Value = 4099.769
$ws.Range("K122").Value = 6719.375100000001
$ws.Range("L122").Value = 12299.307
$ws.Range("M122").Value = -4269.375100000001
$ws.Range("N122").Value = -17199.307

$ws.Range("H132").Value = 3039.9636
$ws.Range("I132").Value = 3004.535
$ws.Range("K132").Value = 9013.605
$ws.Range("M132").Value = -6483.605

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1749
$ws.Range("I132").Value = 1835.5769
$ws.Range("J132").Value = 1373.8334
$ws.Range("K132").Value = 5506.7307
$ws.Range("L132").Value = 4121.5002
$ws.Range("M132").Value = -2976.7307
$ws.Range("N132").Value = -9181.5002
